$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Periodo Mora" (column E) values were listed descending (2103 .. 2003);
# update the database so the list of periods now reads ascending (2003 .. 2103).
$periods = @("2003","2004","2005","2006","2007","2008","2009","2010","2011","2012","2101","2102","2103")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# "Valor Mora" (column F): the lower value (26919) travels together with
# period 2103, which is now the last row (28) instead of the first (16).
for ($row = 16; $row -le 28; $row++) {
    $ws.Range("F$row").Value = 35112
}
$ws.Range("F28").Value = 26919
